$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns keep their original text formatting
# (so values like "1.00" or "0.110" are not auto-converted to numbers)
$ws.Range('D2:E51').NumberFormat = '@'

$ws.Range('D2').Value = '29.684.24'
$ws.Range('E2').Value = '  +1.08%  '
$ws.Range('D3').Value = '1.612.16'
$ws.Range('E3').Value = '  +0.75%  '
$ws.Range('D4').Value = '0.993'
$ws.Range('E4').Value = '  -0.70%  '
$ws.Range('D5').Value = '212.43'
$ws.Range('E5').Value = '  +0.22%  '
$ws.Range('E6').Value = '  +0.59%  '
$ws.Range('D7').Value = '0.992'
$ws.Range('E7').Value = '  -0.68%  '
$ws.Range('D8').Value = '29.04'
$ws.Range('E8').Value = '  +9.99%  '
$ws.Range('E9').Value = '  +3.07%  '
$ws.Range('E10').Value = '  +1.74%  '
$ws.Range('D11').Value = '0.0909'
$ws.Range('E11').Value = '  -0.03%  '
$ws.Range('D12').Value = '1.843.86'
$ws.Range('E12').Value = '  +0.65%  '
$ws.Range('D13').Value = '1.616.86'
$ws.Range('E13').Value = '  +0.71%  '
$ws.Range('D14').Value = '0.568'
$ws.Range('E14').Value = '  +6.90%  '
$ws.Range('D15').Value = '3.87'
$ws.Range('E15').Value = '  +5.11%  '
$ws.Range('D16').Value = '29.664.33'
$ws.Range('E16').Value = '  +0.74%  '
$ws.Range('D17').Value = '8.88'
$ws.Range('E17').Value = '  +16.68%  '
$ws.Range('E18').Value = '  +1.51%  '
$ws.Range('D19').Value = '240.92'
$ws.Range('E19').Value = '  +0.71%  '
$ws.Range('D20').Value = '0.0₃0709'
$ws.Range('E20').Value = '  +3.21%  '
$ws.Range('E21').Value = '  -0.46%  '
$ws.Range('E22').Value = '  +2.93%  '
$ws.Range('D23').Value = '9.65'
$ws.Range('E23').Value = '  +5.96%  '
$ws.Range('E24').Value = '  +1.25%  '
$ws.Range('D25').Value = '156.18'
$ws.Range('E25').Value = '  +1.03%  '
$ws.Range('D26').Value = '15.61'
$ws.Range('E26').Value = '  +2.55%  '
$ws.Range('D27').Value = '0.110'
$ws.Range('D28').Value = '6.58'
$ws.Range('E28').Value = '  +3.62%  '
$ws.Range('E29').Value = '  -0.69%  '
$ws.Range('E30').Value = '  +3.31%  '
$ws.Range('E31').Value = '  +2.28%  '
$ws.Range('E32').Value = '  +3.08%  '
$ws.Range('E33').Value = '  +3.67%  '
$ws.Range('D34').Value = '1.423.31'
$ws.Range('E34').Value = '  +0.97%  '
$ws.Range('D35').Value = '1.63'
$ws.Range('E35').Value = '  +7.01%  '
$ws.Range('E36').Value = '  +1.48%  '
$ws.Range('E37').Value = '  +2.31%  '
$ws.Range('D38').Value = '2.27'
$ws.Range('E38').Value = '  -1.27%  '
$ws.Range('E39').Value = '  +3.14%  '
$ws.Range('D40').Value = '0.555'
$ws.Range('E40').Value = '  +3.89%  '
$ws.Range('D41').Value = '0.0504'
$ws.Range('E41').Value = '  +3.87%  '
$ws.Range('D42').Value = '0.828'
$ws.Range('E42').Value = '  +4.60%  '
$ws.Range('E43').Value = '  +0.36%  '
$ws.Range('D44').Value = '69.86'
$ws.Range('E44').Value = '  +6.76%  '
$ws.Range('D45').Value = '53.66'
$ws.Range('E45').Value = '  +3.22%  '
$ws.Range('D46').Value = '0.992'
$ws.Range('E46').Value = '  -0.64%  '
$ws.Range('D47').Value = '1.00'
$ws.Range('E47').Value = '  +17.63%  '
$ws.Range('D48').Value = '5.43'
$ws.Range('D49').Value = '1.752.65'
$ws.Range('E49').Value = '  +0.55%  '
$ws.Range('D50').Value = '87.69'
$ws.Range('E50').Value = '  +1.41%  '
$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').Value = '0.0₆0108'
$ws.Range('E51').Value = '  +6.44%  '
